# "Working with 2 enemies"
#
# Append two new plain paragraphs (no list/numbering formatting) to the
# end of the document body, right after the existing
# "If vision of the player is lost ..." bullet:
#
#   Official Bug List:
#   Enemy starts shooting through all at bottom right corner of carbon
#   room (possibly not included in room assign and is set to outside?)
#
# Word's Range.Text-based insertion (InsertAfter / InsertParagraphAfter)
# always inherits the paragraph/list formatting of the insertion point,
# which would give the new paragraphs the "ListParagraph" numbered-list
# style used by the preceding bullets. The target paragraphs are plain
# (no <w:pPr> at all), matching the earlier "Enemy Logic:" heading, so
# we build them with Range.InsertXML instead — it lets us hand Word a
# literal WordprocessingML fragment (no pPr => default/plain paragraph)
# to splice in at the end of the story, rather than cloning whatever
# formatting sits at the insertion point.

$d = $word.ActiveDocument

# Collapsed range positioned at the very end of the document's main
# text story (just before the final paragraph mark / sectPr).
$endRange = $d.Range($d.Content.End, $d.Content.End)

$bugListHeading = "Official Bug List:"
$bugListItem = "Enemy starts shooting through all at bottom right corner of carbon room (possibly not included in room assign and is set to outside?)"

$xmlFragment = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t>$bugListHeading</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>$bugListItem</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$endRange.InsertXML($xmlFragment)
